$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the birthday date for the first entry (B1): 08.02 -> 09.02
$ws.Range("B1").Value = "09.02"

# Update the selected cell shown in the sheet view
$ws.Range("C10").Select()
